# Revert "FINFLUX-2918 Nabkisan 3months compunding and subsidy scenarios"
# Restores the pre-change numeric values (and the couple of cells whose
# number format reverts from 2-decimal "#,##0.00" back to whole "#,##0")
# across the Summary, Repayment schedule and Transactions sheets, plus the
# on-screen selections that were recorded at save time.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B2").Value = 1634.47
$wsSummary.Range("E2").Value = 8418.86
$wsSummary.Range("F2").Value = 6702.95

$wsSummary.Range("A3").Value = 356.67
$wsSummary.Range("E3").Value = 256.14
$wsSummary.Range("F3").Value = 237.05

# ---------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

$wsRepay.Range("F3").Value = 1634.47
$wsRepay.Range("G3").Value = 8418.86
$wsRepay.Range("K3").NumberFormat = "#,##0"
$wsRepay.Range("K3").Value = 1735
$wsRepay.Range("L3").NumberFormat = "#,##0"
$wsRepay.Range("L3").Value = 1735

$wsRepay.Range("F4").Value = 1650.81
$wsRepay.Range("G4").Value = 6768.05
$wsRepay.Range("H4").Value = 84.19

$wsRepay.Range("F5").Value = 1667.32
$wsRepay.Range("G5").Value = 5100.73
$wsRepay.Range("H5").Value = 67.68

$wsRepay.Range("F6").Value = 1683.99
$wsRepay.Range("G6").Value = 3416.74
$wsRepay.Range("H6").Value = 51.01

$wsRepay.Range("F7").Value = 1700.83
$wsRepay.Range("G7").Value = 1715.91
$wsRepay.Range("H7").Value = 34.17

$wsRepay.Range("F8").Value = 1715.91
$wsRepay.Range("H8").Value = 19.09

# ---------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")

$wsTrans.Range("A2").Value = 117
$wsTrans.Range("E2").NumberFormat = "#,##0"
$wsTrans.Range("E2").Value = 1735
$wsTrans.Range("F2").Value = 1634.47
$wsTrans.Range("J2").Value = 8418.86

$wsTrans.Range("A3").Value = 116
$wsTrans.Range("A4").Value = 115
$wsTrans.Range("A5").Value = 114

# ---------------------------------------------------------------------
# Restore the saved selections / scroll state on each sheet. Doing the
# "Transactions" sheet last leaves it as the active / tab-selected sheet,
# matching the workbook's recorded activeTab.
# ---------------------------------------------------------------------
$wsSummary.Range("C9").Select()

$wsRepay.Range("K11").Select()

$wsTrans.Range("C8").Select()
